$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two placeholder rows ("* Добавить этаж к текущей схеме" / "*Удалить")
$ws.Rows("18:19").Delete()

# Highlight the (now renumbered) rows 18-20 with an orange fill and
# write the new "ВЫЧЕРКНУТО" (crossed out) column E
$ws.Range("A18:D20").Interior.Color = 49407
$ws.Range("B18:C20").HorizontalAlignment = -4108

$ws.Range("E18").Value = "ВЫЧЕРКНУТО"
$ws.Range("E19").Value = "ВЫЧЕРКНУТО"
$ws.Range("E20").Value = "ВЫЧЕРКНУТО"
$ws.Range("E18:E20").Interior.Color = 49407

$ws.Columns("E").ColumnWidth = 13.140625

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("F20").Select()
